$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 319
$ws.Range("A319").Value = 318
$ws.Range("B319").Value = 'Monday, Jan 16'
$ws.Range("C319").Value = '5:00 AM'
$ws.Range("D319").Value = 'UNKNOWN'
$ws.Range("E319").Value = 'Gdansk'
$ws.Range("F319").Value = '(GDN)'
$ws.Range("G319").Value = 'Enter Air '
$ws.Range("H319").Value = 'B738'
$ws.Range("I319").Value = '(SP-ESD)'
$ws.Range("J319").Value = '5:20 AM'
$ws.Range("L319").Value = '0 hours, 20 minutes'

# Row 320
$ws.Range("A320").Value = 319
$ws.Range("B320").Value = 'Monday, Jan 16'
$ws.Range("C320").Value = '5:00 AM'
$ws.Range("D320").Value = 'UNKNOWN'
$ws.Range("E320").Value = 'Dammam'
$ws.Range("F320").Value = '(DMM)'
$ws.Range("G320").Value = 'Smartwings '
$ws.Range("H320").Value = 'B38M'
$ws.Range("I320").Value = '(OK-SWC)'
$ws.Range("J320").Value = '5:18 AM'
$ws.Range("L320").Value = '0 hours, 18 minutes'

# Row 321
$ws.Range("A321").Value = 320
$ws.Range("B321").Value = 'Monday, Jan 16'
$ws.Range("C321").Value = '5:10 AM'
$ws.Range("D321").Value = 'BO625'
$ws.Range("E321").Value = 'Madrid'
$ws.Range("F321").Value = '(MAD)'
$ws.Range("G321").Value = 'Bluebird Nordic '
$ws.Range("H321").Value = 'B734'
$ws.Range("I321").Value = '(TF-BBN)'
$ws.Range("J321").Value = '5:52 AM'
$ws.Range("L321").Value = '0 hours, 42 minutes'

# Row 322
$ws.Range("A322").Value = 321
$ws.Range("B322").Value = 'Monday, Jan 16'
$ws.Range("C322").Value = '5:40 AM'
$ws.Range("D322").Value = 'LO3880'
$ws.Range("E322").Value = 'Warsaw'
$ws.Range("F322").Value = '(WAW)'
$ws.Range("G322").Value = 'LOT '
$ws.Range("H322").Value = 'E190'
$ws.Range("I322").Value = '(SP-LMG)'
$ws.Range("J322").Value = '6:05 AM'
$ws.Range("L322").Value = '0 hours, 25 minutes'

# Row 323
$ws.Range("A323").Value = 322
$ws.Range("B323").Value = 'Monday, Jan 16'
$ws.Range("C323").Value = '6:15 AM'
$ws.Range("D323").Value = 'W61001'
$ws.Range("E323").Value = 'London'
$ws.Range("F323").Value = '(LTN)'
$ws.Range("G323").Value = 'Wizz Air '
$ws.Range("H323").Value = 'A21N'
$ws.Range("I323").Value = '(HA-LZD)'
$ws.Range("J323").Value = '7:36 AM'
$ws.Range("L323").Value = '1 hours, 21 minutes'

# Row 324
$ws.Range("A324").Value = 323
$ws.Range("B324").Value = 'Monday, Jan 16'
$ws.Range("C324").Value = '6:20 AM'
$ws.Range("D324").Value = 'W61043'
$ws.Range("E324").Value = 'Catania'
$ws.Range("F324").Value = '(CTA)'
$ws.Range("G324").Value = 'Wizz Air '
$ws.Range("H324").Value = 'A21N'
$ws.Range("I324").Value = '(HA-LZJ)'
$ws.Range("J324").Value = '6:34 AM'
$ws.Range("L324").Value = '0 hours, 14 minutes'

# Row 325
$ws.Range("A325").Value = 324
$ws.Range("B325").Value = 'Monday, Jan 16'
$ws.Range("C325").Value = '6:30 AM'
$ws.Range("D325").Value = 'W61091'
$ws.Range("E325").Value = 'Dortmund'
$ws.Range("F325").Value = '(DTM)'
$ws.Range("G325").Value = 'Wizz Air '
$ws.Range("H325").Value = 'A321'
$ws.Range("I325").Value = '(HA-LTC)'
$ws.Range("J325").Value = '6:47 AM'
$ws.Range("L325").Value = '0 hours, 17 minutes'

# Row 326
$ws.Range("A326").Value = 325
$ws.Range("B326").Value = 'Monday, Jan 16'
$ws.Range("C326").Value = '6:40 AM'
$ws.Range("D326").Value = 'W61163'
$ws.Range("E326").Value = 'Malta'
$ws.Range("F326").Value = '(MLA)'
$ws.Range("G326").Value = 'Wizz Air '
$ws.Range("H326").Value = 'A321'
$ws.Range("I326").Value = '(HA-LXN)'
$ws.Range("J326").Value = '7:16 AM'
$ws.Range("L326").Value = '0 hours, 36 minutes'

# Row 327
$ws.Range("A327").Value = 326
$ws.Range("B327").Value = 'Monday, Jan 16'
$ws.Range("C327").Value = '6:45 AM'
$ws.Range("D327").Value = 'LH1363'
$ws.Range("E327").Value = 'Frankfurt'
$ws.Range("F327").Value = '(FRA)'
$ws.Range("G327").Value = 'Lufthansa '
$ws.Range("H327").Value = 'CRJ9'
$ws.Range("I327").Value = '(D-ACNO)'
$ws.Range("J327").Value = '7:05 AM'
$ws.Range("L327").Value = '0 hours, 20 minutes'

# Row 328
$ws.Range("A328").Value = 327
$ws.Range("B328").Value = 'Monday, Jan 16'
$ws.Range("C328").Value = '6:45 AM'
$ws.Range("D328").Value = 'UNKNOWN'
$ws.Range("E328").Value = 'Geilenkirchen'
$ws.Range("F328").Value = '(GKE)'
$ws.Range("G328").Value = 'Enter Air '
$ws.Range("H328").Value = 'B738'
$ws.Range("I328").Value = '(SP-ESH)'
$ws.Range("J328").Value = '7:02 AM'
$ws.Range("L328").Value = '0 hours, 17 minutes'

# Row 329
$ws.Range("A329").Value = 328
$ws.Range("B329").Value = 'Monday, Jan 16'
$ws.Range("C329").Value = '7:00 AM'
$ws.Range("D329").Value = 'FR6391'
$ws.Range("E329").Value = 'London'
$ws.Range("F329").Value = '(STN)'
$ws.Range("G329").Value = 'Ryanair '
$ws.Range("H329").Value = 'B738'
$ws.Range("I329").Value = '(SP-RSB)'
$ws.Range("J329").Value = '7:27 AM'
$ws.Range("L329").Value = '0 hours, 27 minutes'

# Row 330
$ws.Range("A330").Value = 329
$ws.Range("B330").Value = 'Monday, Jan 16'
$ws.Range("C330").Value = '7:35 AM'
$ws.Range("D330").Value = 'FR6892'
$ws.Range("E330").Value = 'Dortmund'
$ws.Range("F330").Value = '(DTM)'
$ws.Range("G330").Value = 'Ryanair '
$ws.Range("H330").Value = 'B738'
$ws.Range("I330").Value = '(SP-RSN)'
$ws.Range("J330").Value = '7:46 AM'
$ws.Range("L330").Value = '0 hours, 11 minutes'

# Row 331
$ws.Range("A331").Value = 330
$ws.Range("B331").Value = 'Monday, Jan 16'
$ws.Range("C331").Value = '7:40 AM'
$ws.Range("D331").Value = 'E47041'
$ws.Range("E331").Value = 'Hurghada'
$ws.Range("F331").Value = '(HRG)'
$ws.Range("G331").Value = 'Enter Air '
$ws.Range("H331").Value = 'B738'
$ws.Range("I331").Value = '(SP-ESI)'
$ws.Range("J331").Value = '7:59 AM'
$ws.Range("L331").Value = '0 hours, 19 minutes'

# Row 332
$ws.Range("A332").Value = 331
$ws.Range("B332").Value = 'Monday, Jan 16'
$ws.Range("C332").Value = '8:30 AM'
$ws.Range("D332").Value = '3Z7312'
$ws.Range("E332").Value = 'Fuerteventura'
$ws.Range("F332").Value = '(FUE)'
$ws.Range("G332").Value = 'Smartwings '
$ws.Range("H332").Value = 'B738'
$ws.Range("I332").Value = '(OK-TSF)'
$ws.Range("J332").Value = '8:48 AM'
$ws.Range("L332").Value = '0 hours, 18 minutes'

# Row 333
$ws.Range("A333").Value = 332
$ws.Range("B333").Value = 'Monday, Jan 16'
$ws.Range("C333").Value = '9:35 AM'
$ws.Range("D333").Value = 'BO951'
$ws.Range("E333").Value = 'Paris'
$ws.Range("F333").Value = '(CDG)'
$ws.Range("G333").Value = 'Bluebird Nordic '
$ws.Range("H333").Value = 'B734'
$ws.Range("I333").Value = '(TF-BBO)'
$ws.Range("J333").Value = '9:30 AM'
$ws.Range("L333").Value = '0 hours, -5 minutes'

# Row 334
$ws.Range("A334").Value = 333
$ws.Range("B334").Value = 'Monday, Jan 16'
$ws.Range("C334").Value = '9:45 AM'
$ws.Range("D334").Value = 'UNKNOWN'
$ws.Range("E334").Value = 'Fuerteventura'
$ws.Range("F334").Value = '(FUE)'
$ws.Range("G334").Value = 'Enter Air '
$ws.Range("H334").Value = 'B738'
$ws.Range("I334").Value = '(SP-ENP)'
$ws.Range("J334").Value = '9:47 AM'
$ws.Range("L334").Value = '0 hours, 2 minutes'
